$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the surviving data rows (2-6) with the new values from the diff.
$ws.Range("A2").Value = 3
$ws.Range("B2").Value = 248

$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 241

$ws.Range("A4").Value = 2
$ws.Range("B4").Value = 183

$ws.Range("A5").Value = 0
$ws.Range("B5").Value = 183

$ws.Range("A6").Value = 4
$ws.Range("B6").Value = 147

# Remove rows 7-11 entirely (they are no longer part of the data set).
$ws.Range("A7:B11").EntireRow.Delete()
